# Generate Report for Handoff
#
# The localization-status workbook tracks two source files:
#   73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md
#   f8d2db94-eebd-4206-a6da-1cf04127c58a.md
#
# A new handoff cycle happened for 73bb9739-...: its status moved from
# "Handed back: in sync with en-US" to "Ready for handoff" with a fresh
# handoff timestamp, and its row now sorts after f8d2db94-... (whose data
# is unchanged) on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value() = "f8d2db94-eebd-4206-a6da-1cf04127c58a.md"
$overview.Range("B2").Value() = "Handed back: in sync with en-US"
$overview.Range("C2").Value() = "Handed back: in sync with en-US"
$overview.Range("D2").Value() = "2016-33-18 12:33:11"

$overview.Range("A3").Value() = "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md"
$overview.Range("B3").Value() = "Ready for handoff"
$overview.Range("C3").Value() = "Ready for handoff"
$overview.Range("D3").Value() = "2016-34-18 12:34:07"

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status |
#   Latest Handoff File | Latest Handoff Datetime | Latest Target File |
#   Latest Handback File | Latest Handback DateTime | Handoff Reason |
#   Dependency From | Error Detail
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Value() = "f8d2db94-eebd-4206-a6da-1cf04127c58a.md"
$zhcn.Range("B2").Value() = ".md"
$zhcn.Range("C2").Value() = "Handed back: in sync with en-US"
$zhcn.Range("D2").Value() = "f8d2db94-eebd-4206-a6da-1cf04127c58a.b4ca52014fd439e0f863e2ffa23315fcf9010504.zh-cn.xlf"
$zhcn.Range("E2").Value() = "2016-03-18 12:33:08"
$zhcn.Range("F2").Value() = "f8d2db94-eebd-4206-a6da-1cf04127c58a.md"
$zhcn.Range("G2").Value() = "f8d2db94-eebd-4206-a6da-1cf04127c58a.b4ca52014fd439e0f863e2ffa23315fcf9010504.zh-cn.xlf"
$zhcn.Range("H2").Value() = "2016-03-18 12:33:36"
$zhcn.Range("I2").Value() = "Include"

$zhcn.Range("A3").Value() = "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md"
$zhcn.Range("B3").Value() = ".md"
$zhcn.Range("C3").Value() = "Ready for handoff"
$zhcn.Range("D3").Value() = "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.11ed010e27266a27c332c8fdd55168dead1be9f6.zh-cn.xlf"
$zhcn.Range("E3").Value() = "2016-03-18 12:34:04"
$zhcn.Range("F3").Value() = "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md"
$zhcn.Range("G3").Value() = "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.11ed010e27266a27c332c8fdd55168dead1be9f6.zh-cn.xlf"
$zhcn.Range("H3").Value() = "2016-03-18 12:33:36"
$zhcn.Range("I3").Value() = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de": same column layout as zh-cn
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Value() = "f8d2db94-eebd-4206-a6da-1cf04127c58a.md"
$dede.Range("B2").Value() = ".md"
$dede.Range("C2").Value() = "Handed back: in sync with en-US"
$dede.Range("D2").Value() = "f8d2db94-eebd-4206-a6da-1cf04127c58a.b4ca52014fd439e0f863e2ffa23315fcf9010504.de-de.xlf"
$dede.Range("E2").Value() = "2016-03-18 12:33:11"
$dede.Range("F2").Value() = "f8d2db94-eebd-4206-a6da-1cf04127c58a.md"
$dede.Range("G2").Value() = "f8d2db94-eebd-4206-a6da-1cf04127c58a.b4ca52014fd439e0f863e2ffa23315fcf9010504.de-de.xlf"
$dede.Range("H2").Value() = "2016-03-18 12:33:43"
$dede.Range("I2").Value() = "Include"

$dede.Range("A3").Value() = "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md"
$dede.Range("B3").Value() = ".md"
$dede.Range("C3").Value() = "Ready for handoff"
$dede.Range("D3").Value() = "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.11ed010e27266a27c332c8fdd55168dead1be9f6.de-de.xlf"
$dede.Range("E3").Value() = "2016-03-18 12:34:07"
$dede.Range("F3").Value() = "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.md"
$dede.Range("G3").Value() = "73bb9739-d02c-4df7-9f3e-d5f82fd5f619.11ed010e27266a27c332c8fdd55168dead1be9f6.de-de.xlf"
$dede.Range("H3").Value() = "2016-03-18 12:33:43"
$dede.Range("I3").Value() = "Include"

